# Convert distributions/treatment parser template:
#  - rename the sheet "Sheet1" -> "Template"
#  - move the header row from row 1 down to row 3
#  - rework the header columns: drop "Len (mm)", add "Collection" and
#    "Len (cm)", and reorder the trailing columns accordingly
#  - size a few of the new/changed columns
#  - leave the selection on H3 (Truck)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet.
$ws.Name = "Template"

# Push the existing header (currently row 1) down to row 3.
$ws.Rows("1:2").Insert()

# Re-write the header row (row 3) in its new column order.
# NB: "Len (cm)" is written before "Collection" so that the two brand-new
# shared strings land in the same append order as the target workbook.
$ws.Range("A3").Value = "Year"
$ws.Range("B3").Value = "Month"
$ws.Range("C3").Value = "Day"
$ws.Range("D3").Value = "Time"
$ws.Range("E3").Value = "Site"
$ws.Range("F3").Value = "Lat"
$ws.Range("G3").Value = "Long"
$ws.Range("H3").Value = "Truck"
$ws.Range("I3").Value = "Driver"
$ws.Range("J3").Value = "Crew"
$ws.Range("K3").Value = "Stock"
$ws.Range("M3").Value = "Tank"
$ws.Range("N3").Value = "NFish"
$ws.Range("O3").Value = "Len (cm)"
$ws.Range("L3").Value = "Collection"
$ws.Range("P3").Value = "Temp"
$ws.Range("Q3").Value = "Comments"

# Drop the old trailing blank column (R) now that the used range ends at Q.
$ws.Columns("R").Delete() | Out-Null

# Widen the new/changed columns (Collection, Len (cm), Comments).
$ws.Columns("L").ColumnWidth = 9.166666666666666
$ws.Columns("O").ColumnWidth = 9.666666666666666
$ws.Columns("Q").ColumnWidth = 11

# Match the saved selection.
$ws.Range("H3").Select() | Out-Null
